$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# Keep the Week_Start_Date column as plain text (matches source file which
# stores these as inline strings, not real Excel dates).
$wsForecast.Range("B2:B17").NumberFormat = "@"

# --- Forecast Comparison sheet: shift Week_Start_Date values forward by one week ---
$wsForecast.Range("B2").Value  = "2025-02-02"
$wsForecast.Range("B3").Value  = "2025-02-09"
$wsForecast.Range("B4").Value  = "2025-02-16"
$wsForecast.Range("B5").Value  = "2025-02-23"
$wsForecast.Range("B6").Value  = "2025-03-02"
$wsForecast.Range("B7").Value  = "2025-03-09"
$wsForecast.Range("B8").Value  = "2025-03-16"
$wsForecast.Range("B9").Value  = "2025-03-23"
$wsForecast.Range("B10").Value = "2025-03-30"
$wsForecast.Range("B11").Value = "2025-04-06"
$wsForecast.Range("B12").Value = "2025-04-13"

$wsForecast.Range("B13").Value = "2025-04-20"
$wsForecast.Range("E13").Value = 1
$wsForecast.Range("F13").Value = 1
$wsForecast.Range("G13").Value = 1
$wsForecast.Range("H13").Value = 2

$wsForecast.Range("B14").Value = "2025-04-27"
$wsForecast.Range("E14").Value = 1
$wsForecast.Range("F14").Value = 1
$wsForecast.Range("G14").Value = 1
$wsForecast.Range("H14").Value = 2

$wsForecast.Range("B15").Value = "2025-05-04"
$wsForecast.Range("G15").Value = 1
$wsForecast.Range("H15").Value = 1

$wsForecast.Range("B16").Value = "2025-05-11"
$wsForecast.Range("G16").Value = 1
$wsForecast.Range("H16").Value = 1

$wsForecast.Range("B17").Value = "2025-05-18"
$wsForecast.Range("E17").Value = 1
$wsForecast.Range("F17").Value = 1
$wsForecast.Range("G17").Value = 1
$wsForecast.Range("H17").Value = 3

# --- Summary sheet updates ---
$wsSummary.Range("B2").NumberFormat  = "@"
$wsSummary.Range("B13").NumberFormat = "@"
$wsSummary.Range("B15").NumberFormat = "@"

$wsSummary.Range("B2").Value  = "2022-12-25 to 2025-01-26"
$wsSummary.Range("B13").Value = "2025-04-20"
$wsSummary.Range("B15").Value = "2025-02-02"
